$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and swap Aptos / InternetComputer rows)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.534.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5056"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3955"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09762"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.543"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.908.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.567"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001139"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06647"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.600.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.278"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.766"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +16.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.127.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.41%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.738"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.639"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.866"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06797"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02445"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.270"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.95%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.73"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.87%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.110"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6437"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.190"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.18%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6098"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.282"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.669"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.07%  "
